$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph "Database: 0.18 (Git tag: cen_cruise_db_v0.18)" and
# the paragraph that starts with "Description:" (the second occurrence,
# the one that still carries the stray "_GoBack" bookmark).
# ---------------------------------------------------------------------------
$verPara = $null
$descPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Database: 0.18 (Git tag: cen_cruise_db_v0.18)`r") {
        $verPara = $p
    }
    if ($t -like "Description:*internal report to provide*") {
        $descPara = $p
    }
}

if ($verPara -eq $null) { throw "version paragraph not found" }
if ($descPara -eq $null) { throw "description paragraph not found" }

$vStart = $verPara.Range.Start

# ---------------------------------------------------------------------------
# 1) "Database: 0." + "1" + "8"  ->  change the "8" run to "9"
#    layout (0-based offsets from paragraph start):
#      0..11  "Database: 0."
#      12     "1"
#      13     "8"   <- change to "9"
# ---------------------------------------------------------------------------
$firstEight = $d.Range($vStart + 13, $vStart + 14)
if ($firstEight.Text -ne "8") { throw "unexpected char at first-eight position: [$($firstEight.Text)]" }
$firstEight.Text = "9"

# ---------------------------------------------------------------------------
# 2) Insert a (new) "_GoBack" bookmark right after that run, mirroring what
#    a live edit in Word leaves behind at the last edited location.
# ---------------------------------------------------------------------------
$gobackSpot = $d.Range($vStart + 14, $vStart + 14)
if ($d.Bookmarks.Exists("_GoBack")) { $d.Bookmarks("_GoBack").Delete() }
$d.Bookmarks.Add("_GoBack", $gobackSpot)

# ---------------------------------------------------------------------------
# 3) Second "8" -> "9" :  "... tag: cen_cruise_db_v0." + "1" + "8" + ")"
# ---------------------------------------------------------------------------
$full = $verPara.Range.Text
$idx = $full.IndexOf("v0.18")
if ($idx -lt 0) { throw "git tag version text not found" }
$secondEight = $d.Range($vStart + $idx + 4, $vStart + $idx + 5)
if ($secondEight.Text -ne "8") { throw "unexpected char at second-eight position: [$($secondEight.Text)]" }
$secondEight.Text = "9"

Write-Host "Version paragraph now: [" $verPara.Range.Text "]"

# ---------------------------------------------------------------------------
# 4) Tidy up the "Description:" paragraph:
#      - remove the stray "_GoBack" bookmark sitting inside the middle of it
#      - collapse the run split so the text reads as two runs:
#          "Description: "  (keeps the lastRenderedPageBreak marker)
#          "This standard ... desired"
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) { $d.Bookmarks("_GoBack").Delete() }

$dStart = $descPara.Range.Start
$full2 = $descPara.Range.Text
$idx2 = $full2.IndexOf("Description:")
$spacePos = $dStart + $idx2 + ("Description:").Length

$spaceRange = $d.Range($spacePos, $spacePos + 1)
if ($spaceRange.Text -ne " ") { throw "unexpected char at space position: [$($spaceRange.Text)]" }
# Round-trip through a placeholder value: an unchanged write is a no-op in
# this runtime, so briefly swap the character to force the surrounding runs
# to coalesce, then restore it.
$spaceRange.Text = "#"
$spaceRange2 = $d.Range($spacePos, $spacePos + 1)
$spaceRange2.Text = " "

Write-Host "Description paragraph now: [" $descPara.Range.Text "]"
